$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of rows 32 and 33 ---
# (Row 32 becomes the old row-33 record; row 33 becomes the old row-32 record.
#  The columns that are identical between the two rows are left untouched.)

# Row 32 -> old row 33's values
$ws.Range("A32").Value = 112213305
$ws.Range("B32").Value = 89517
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 5447
$ws.Range("F32").Value = "Vedticka"
$ws.Range("G32").Value = "Fuscoporia viticola"
$ws.Range("H32").Value = "(Schwein.) Murrill"
$ws.Range("Q32").Value = 515748
$ws.Range("R32").Value = 6704727

# Row 33 -> old row 32's values
$ws.Range("A33").Value = 112213272
$ws.Range("B33").Value = 89553
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q33").Value = 515738
$ws.Range("R33").Value = 6704726

# --- Append new row 35 ---
$ws.Range("A35").Value = 112539791
$ws.Range("B35").Value = 90858
$ws.Range("C35").Value = "Ovaliderad"
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 5449
$ws.Range("F35").Value = "Svart taggsvamp"
$ws.Range("G35").Value = "Phellodon niger"
$ws.Range("H35").Value = "(Fr.:Fr.) P.Karst."

# I35 must stay a text "1" (like the rest of the "Antal" column), not a number.
$ws.Range("I35").Value = "'1"
$ws.Range("I35").Style = "Normal"

$ws.Range("J35").Value = "fruktkroppar"

# K35 / N35 are present but empty text cells.
$ws.Range("K35").Value = "'"
$ws.Range("K35").Style = "Normal"
$ws.Range("N35").Value = "'"
$ws.Range("N35").Style = "Normal"

$ws.Range("P35").Value = "Simsbodarna, Dlr"
$ws.Range("Q35").Value = 515374
$ws.Range("R35").Value = 6704951
$ws.Range("S35").Value = 1
$ws.Range("T35").Value = "Dalarna"
$ws.Range("U35").Value = "Borlänge"
$ws.Range("V35").Value = "Dalarna"
$ws.Range("W35").Value = "Stora Tuna"

# Y35 / AA35 are text dates (e.g. "2023-10-05"), not real Excel date serials.
$ws.Range("Y35").Value = "'2023-10-05"
$ws.Range("Y35").Style = "Normal"
$ws.Range("Z35").Value = "16:30"
$ws.Range("AA35").Value = "'2023-10-05"
$ws.Range("AA35").Style = "Normal"
$ws.Range("AB35").Value = "16:30"

$ws.Range("AC35").Value = "En mindre fk."
$ws.Range("AD35").Value = $false
$ws.Range("AE35").Value = $false

$ws.Range("AF35").Value = "'"
$ws.Range("AF35").Style = "Normal"

$ws.Range("AG35").Value = $false

$ws.Range("AT35").Value = "'"
$ws.Range("AT35").Style = "Normal"

$ws.Range("AW35").Value = "Lars-Erik Nilsson"
$ws.Range("AX35").Value = "Lars-Erik Nilsson, Bo karlstens, Anna-Lena Thommson, Niklas Rehn, Matilda Elgerud, Malte Lindberg, Holger Martinussen"

$ws.Range("AY35").Value = "'"
$ws.Range("AY35").Style = "Normal"
